# GD - Justificacion de anulacion
# Resumen sheet: insert 2 columns before old column F ("Creacion"), relabel the
# date columns, rename "Materia" -> moves with the new date-style format, and
# "referencia" -> "Referencia". Also widen the AutoFilter ranges on both
# sheets (and keep the _FilterDatabase defined names in sync), matching the
# target workbook state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Resumen (sheet 1)
# ---------------------------------------------------------------------------
$resumen = $wb.Worksheets.Item("Resumen")

# Insert two new columns at F:G (old F/G "Creacion"/"Fecha ingreso" shift to H/I)
$resumen.Range("F1:G1").EntireColumn.Insert()

# New F1 / G1 headers
$resumen.Range("F1").Value = "Inicio proceso"
$resumen.Range("G1").Value = "Término proceso"
$resumen.Range("F1:G1").Font.Bold = $true
$resumen.Range("F1:G1").NumberFormat = "m/d/yy h:mm"

# Old "Creacion"/"Fecha ingreso" headers (now H1/I1) get renamed
$resumen.Range("H1").Value = "Fecha creación GD"
$resumen.Range("I1").Value = "Fecha ingreso GD"

# "Materia" (now J1) switches from the plain bold style to the date-style bold
# formatting used by its neighbours (content unchanged)
$resumen.Range("J1").Font.Bold = $true
$resumen.Range("J1").NumberFormat = "m/d/yy h:mm"

# "referencia" (now K1) gets capitalised
$resumen.Range("K1").Value = "Referencia"

# New column widths (best effort - original values were derived from Excel's
# AutoFit over data rows that are not present in this workbook)
$resumen.Columns.Item(7).ColumnWidth = 20.33
$resumen.Columns.Item(8).ColumnWidth = 18.88
$resumen.Columns.Item(9).ColumnWidth = 17.88
$resumen.Columns.Item(10).ColumnWidth = 17.88
$resumen.Columns.Item(17).ColumnWidth = 17.45

# Refresh the AutoFilter over the (now wider/taller) data range
if ($resumen.AutoFilterMode) {
    $resumen.AutoFilterMode = $false
}
$resumen.Range("A1:S157").AutoFilter()

# ---------------------------------------------------------------------------
# Detalle (sheet 2)
# ---------------------------------------------------------------------------
$detalle = $wb.Worksheets.Item("Detalle")

if ($detalle.AutoFilterMode) {
    $detalle.AutoFilterMode = $false
}
$detalle.Range("A1:K667").AutoFilter()

# ---------------------------------------------------------------------------
# Keep the workbook-level _FilterDatabase defined names in sync with the new
# AutoFilter ranges (the host does not do this automatically)
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Resumen!_FilterDatabase") {
        $n.RefersTo = "=Resumen!`$A`$1:`$S`$157"
    }
    if ($n.Name -eq "Detalle!_FilterDatabase") {
        $n.RefersTo = "=Detalle!`$A`$1:`$K`$667"
    }
}

Write-Host "edit complete"
